$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose content differs between row 35 and row 36 and must be swapped.
$cols = @("A","B","E","F","G","H","I","L","M","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr35 = "$col" + "35"
    $addr36 = "$col" + "36"
    $v35 = $ws.Range($addr35).Value2
    $v36 = $ws.Range($addr36).Value2

    # "Antal" (I) holds small counts stored as plain TEXT in the source data
    # (e.g. "2"), not numbers. A bare numeric-looking string would be
    # auto-coerced to a real number by normal value assignment, so force it
    # back to text with a leading apostrophe and then strip the resulting
    # quote-prefix formatting so only the value (not the cell style) changes.
    if ($col -eq "I") {
        if ($v36 -eq $null -or $v36 -eq "") {
            $ws.Range($addr35).Value2 = ""
        } else {
            $ws.Range($addr35).Value = "'" + $v36
            $ws.Range($addr35).Style = "Normal"
        }
        if ($v35 -eq $null -or $v35 -eq "") {
            $ws.Range($addr36).Value2 = ""
        } else {
            $ws.Range($addr36).Value = "'" + $v35
            $ws.Range($addr36).Style = "Normal"
        }
    } else {
        $ws.Range($addr35).Value2 = $v36
        $ws.Range($addr36).Value2 = $v35
    }
}
